$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2845399081707001
$ws.Range("B1").Value = 0.2194613963365555
$ws.Range("C1").Value = 0.1882201135158539
$ws.Range("D1").Value = 0.1990492194890976
$ws.Range("E1").Value = 0.2361022680997849
